# Append a new row 36 (next sensor reading) to each of the four sheets,
# mirroring the structure of the existing rows (e.g. row 35).

$wb = $excel.ActiveWorkbook

$rowsBySheet = @{
    "ROW35-FE-LIFTER" = @{
        A = "2025-03-05 19:42:06"
        B = "0x01,0x90 "
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"
        D = "0x01,0x90,"
        E = "0x d"
        F = 400
        G = "'568631262647113770877196"
        H = 400
        I = 13
    }
    "ROW35-MID-LIFTER" = @{
        A = "2025-03-05 19:29:35"
        B = "0x01,0x90 "
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
        D = "0x01,0x90,"
        E = "0x e"
        F = 400
        G = "'568631262647113770942732"
        H = 400
        I = 14
    }
    "ROW02-FE-LIFTER" = @{
        A = "2025-03-05 19:51:45"
        B = "0x01,0x90 "
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x06,0x41,0x0c,"
        D = "0x01,0x90,"
        E = "0xff"
        F = 400
        G = "'568631262647113769959692"
        H = 400
        I = 255
    }
    "ROW02-MID-LIFTER" = @{
        A = "2025-03-05 19:41:15"
        B = "0x01,0x90 "
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x01,0x90,"
        E = "0x 3"
        F = 400
        G = "'568631262647113769959692"
        H = 400
        I = 3
    }
}

foreach ($sheetName in $rowsBySheet.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $data = $rowsBySheet[$sheetName]
    $newRow = $ws.UsedRange.Rows.Count + 1

    $ws.Cells.Item($newRow, 1).Value2 = $data.A
    $ws.Cells.Item($newRow, 2).Value2 = $data.B
    $ws.Cells.Item($newRow, 3).Value2 = $data.C
    $ws.Cells.Item($newRow, 4).Value2 = $data.D
    $ws.Cells.Item($newRow, 5).Value2 = $data.E
    $ws.Cells.Item($newRow, 6).Value2 = $data.F
    $ws.Cells.Item($newRow, 7).Value2 = $data.G
    $ws.Cells.Item($newRow, 8).Value2 = $data.H
    $ws.Cells.Item($newRow, 9).Value2 = $data.I
}
